# Generate Report for Handback
# Adds a new handback record (ce35765b-08c0-4403-b1ec-e01e4fb0289b) as row 4
# on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$guidMd   = "ce35765b-08c0-4403-b1ec-e01e4fb0289b.md"
$guidPath = "e2e\ce35765b-08c0-4403-b1ec-e01e4fb0289b.md"
$xlfZh    = "ce35765b-08c0-4403-b1ec-e01e4fb0289b.f6b2b13c56839ee9e388bfb0cabd6853fde44323.zh-cn.xlf"
$xlfDe    = "ce35765b-08c0-4403-b1ec-e01e4fb0289b.f6b2b13c56839ee9e388bfb0cabd6853fde44323.de-de.xlf"

$dateHandoff   = "2016-08-12 16:55:51"
$dateZhGen     = "2016-08-12 16:55:44"
$dateZhBack    = "2016-08-12 16:56:17"
$dateDeBack    = "2016-08-12 16:56:26"

$statusInSync  = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $guidMd
$wsOverview.Range("B4").Value = $guidPath
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = $statusInSync
$wsOverview.Range("F4").Value = $statusInSync
$wsOverview.Range("G4").Value = $dateHandoff
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Range("B4").Style = "HyperLink"
$linkOverview = $wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/f6b2b13c56839ee9e388bfb0cabd6853fde44323/e2e/ce35765b-08c0-4403-b1ec-e01e4fb0289b.md", "", "", $guidPath)

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $guidMd
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'True"
$wsZh.Range("G4").Value = $xlfZh
$wsZh.Range("H4").Value = $dateZhGen
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $guidMd
$wsZh.Range("J4").Value = $xlfZh
$wsZh.Range("K4").Value = $dateZhBack
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = "'"
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = "'"

$wsZh.Range("A4").Style = "HyperLink"
$wsZh.Range("I4").Style = "HyperLink"
$linkZhA = $wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/f6b2b13c56839ee9e388bfb0cabd6853fde44323/e2e/ce35765b-08c0-4403-b1ec-e01e4fb0289b.md", "", "", $guidMd)
$linkZhI = $wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f6b2b13c56839ee9e388bfb0cabd6853fde44323/e2e/ce35765b-08c0-4403-b1ec-e01e4fb0289b.md", "", "", $guidMd)

$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $guidMd
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'True"
$wsDe.Range("G4").Value = $xlfDe
$wsDe.Range("H4").Value = $dateHandoff
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $guidMd
$wsDe.Range("J4").Value = $xlfDe
$wsDe.Range("K4").Value = $dateDeBack
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = "'"
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = "'"

$wsDe.Range("A4").Style = "HyperLink"
$wsDe.Range("I4").Style = "HyperLink"
$linkDeA = $wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/f6b2b13c56839ee9e388bfb0cabd6853fde44323/e2e/ce35765b-08c0-4403-b1ec-e01e4fb0289b.md", "", "", $guidMd)
$linkDeI = $wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/f6b2b13c56839ee9e388bfb0cabd6853fde44323/e2e/ce35765b-08c0-4403-b1ec-e01e4fb0289b.md", "", "", $guidMd)

$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.Resize($wsDe.Range("A1:P4"))
